$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header "E" in L1
$ws.Range("L1").Value = "E"

# Fill L2:L156 with 252 for every data row
$ws.Range("L2:L156").Value = 252

# Reflect the final cursor/selection position left by the author after
# adding the new column of data (row 157 is just past the last data row).
$ws.Range("P157").Select()
